$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The test data credentials were rotated as part of the Selenium 4.2 update:
#   A2 (username's paired value) -> "mngr628740"
#   B2 (password's paired value) -> "dApAjar"
$ws.Range("A2").Value = "mngr628740"
$ws.Range("B2").Value = "dApAjar"

# Reflect the final user selection resting on B2 (single cell, not the whole row).
$null = $ws.Range("B2").Select()
